$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:I15").ClearContents()

# Ensure column C uses text format so numeric-looking strings are preserved as text
$ws.Range("C2:C15").NumberFormat = "@"

$ws.Range("A2").Value = "(Intercept)"
$ws.Range("B2").Value = 1.812
$ws.Range("C2").Value = "0.128615316640211"
$ws.Range("D2").Value = 0.1037783132149098
$ws.Range("E2").Value = 105495
$ws.Range("F2").Value = "Model 1d"
$ws.Range("G2").Value = 1.641284674761474
$ws.Range("H2").Value = 1.982715325238527
$ws.Range("I2").Value = 38

$ws.Range("A3").Value = "lrscale"
$ws.Range("B3").Value = -0.239
$ws.Range("C3").Value = "0.124838248068265"
$ws.Range("D3").Value = 0.1035949760193604
$ws.Range("E3").Value = 105495
$ws.Range("F3").Value = "Model 1d"
$ws.Range("G3").Value = -0.4094137355518479
$ws.Range("H3").Value = -0.06858626444815205
$ws.Range("I3").Value = 38

$ws.Range("A4").Value = "legacyLeft-Wing"
$ws.Range("B4").Value = -0.892
$ws.Range("C4").Value = "-0.102597411262751"
$ws.Range("D4").Value = 0.1304675345844536
$ws.Range("E4").Value = 105495
$ws.Range("F4").Value = "Model 1d"
$ws.Range("G4").Value = -1.106619094391426
$ws.Range("H4").Value = -0.6773809056085739
$ws.Range("I4").Value = 38

$ws.Range("A5").Value = "legacyRight-Wing"
$ws.Range("B5").Value = 0.191
$ws.Range("C5").Value = "0.128615316640211"
$ws.Range("D5").Value = 0.1785053454107153
$ws.Range("E5").Value = 105495
$ws.Range("F5").Value = "Model 1d"
$ws.Range("G5").Value = -0.1026412932006267
$ws.Range("H5").Value = 0.4846412932006267
$ws.Range("I5").Value = 38

$ws.Range("A6").Value = "age"
$ws.Range("B6").Value = -0.144
$ws.Range("C6").Value = "0.124838248068265"
$ws.Range("D6").Value = 0.01627118319622274
$ws.Range("E6").Value = 105495
$ws.Range("F6").Value = "Model 1d"
$ws.Range("G6").Value = -0.1707660963577864
$ws.Range("H6").Value = -0.1172339036422136
$ws.Range("I6").Value = 38

$ws.Range("A7").Value = "educ"
$ws.Range("B7").Value = 0.305
$ws.Range("C7").Value = "-0.102597411262751"
$ws.Range("D7").Value = 0.00865563850457194
$ws.Range("E7").Value = 105495
$ws.Range("F7").Value = "Model 1d"
$ws.Range("G7").Value = 0.2907614746599791
$ws.Range("H7").Value = 0.3192385253400208
$ws.Range("I7").Value = 38

$ws.Range("A8").Value = "polint"
$ws.Range("B8").Value = -0.458
$ws.Range("C8").Value = "0.128615316640211"
$ws.Range("D8").Value = 0.01004649655726492
$ws.Range("E8").Value = 105495
$ws.Range("F8").Value = "Model 1d"
$ws.Range("G8").Value = -0.4745264868367008
$ws.Range("H8").Value = -0.4414735131632992
$ws.Range("I8").Value = 38

$ws.Range("A9").Value = "sexMale"
$ws.Range("B9").Value = 0.014
$ws.Range("C9").Value = "0.124838248068265"
$ws.Range("D9").Value = 0.005768771990899003
$ws.Range("E9").Value = 105495
$ws.Range("F9").Value = "Model 1d"
$ws.Range("G9").Value = 0.00451037007497114
$ws.Range("H9").Value = 0.02348962992502886
$ws.Range("I9").Value = 38

$ws.Range("A10").Value = "surveyevs2008"
$ws.Range("B10").Value = -0.106
$ws.Range("C10").Value = "-0.102597411262751"
$ws.Range("D10").Value = 0.008332330897955476
$ws.Range("E10").Value = 105495
$ws.Range("F10").Value = "Model 1d"
$ws.Range("G10").Value = -0.1197066843271368
$ws.Range("H10").Value = -0.09229331567286324
$ws.Range("I10").Value = 38

$ws.Range("A11").Value = "surveywvs1994"
$ws.Range("B11").Value = -0.103
$ws.Range("C11").Value = "0.128615316640211"
$ws.Range("D11").Value = 0.01008306006889319
$ws.Range("E11").Value = 105495
$ws.Range("F11").Value = "Model 1d"
$ws.Range("G11").Value = -0.1195866338133293
$ws.Range("H11").Value = -0.0864133661866707
$ws.Range("I11").Value = 38

$ws.Range("A12").Value = "surveywvs1999"
$ws.Range("B12").Value = 0.18
$ws.Range("C12").Value = "0.124838248068265"
$ws.Range("D12").Value = 0.01654601227902841
$ws.Range("E12").Value = 105495
$ws.Range("F12").Value = "Model 1d"
$ws.Range("G12").Value = 0.1527818098009983
$ws.Range("H12").Value = 0.2072181901990017
$ws.Range("I12").Value = 38

$ws.Range("A13").Value = "surveywvs2005"
$ws.Range("B13").Value = 0.299
$ws.Range("C13").Value = "-0.102597411262751"
$ws.Range("D13").Value = 0.01099273797156308
$ws.Range("E13").Value = 105495
$ws.Range("F13").Value = "Model 1d"
$ws.Range("G13").Value = 0.2809169460367787
$ws.Range("H13").Value = 0.3170830539632212
$ws.Range("I13").Value = 38

$ws.Range("A14").Value = "lrscale:legacyLeft-Wing"
$ws.Range("B14").Value = 0.571
$ws.Range("C14").Value = "0.128615316640211"
$ws.Range("D14").Value = 0.1307651844314123
$ws.Range("E14").Value = 105495
$ws.Range("F14").Value = "Model 1d"
$ws.Range("G14").Value = 0.3558912716103266
$ws.Range("H14").Value = 0.7861087283896733
$ws.Range("I14").Value = 38

$ws.Range("A15").Value = "lrscale:legacyRight-Wing"
$ws.Range("B15").Value = -0.26
$ws.Range("C15").Value = "0.124838248068265"
$ws.Range("D15").Value = 0.1790310012811109
$ws.Range("E15").Value = 105495
$ws.Range("F15").Value = "Model 1d"
$ws.Range("G15").Value = -0.5545059971074274
$ws.Range("H15").Value = 0.03450599710742741
$ws.Range("I15").Value = 38

# Reset style on column C cells so no stray unused style persists (keep default style)
$ws.Range("C2:C15").Style = "Normal"
